$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(225, 1).Value = 4
$ws.Cells.Item(225, 2).Value = 19
$ws.Cells.Item(225, 3).Value = 5
$ws.Cells.Item(225, 4).Value = 1
$ws.Cells.Item(225, 5).Formula = "=B225+D225"

$ws.Cells.Item(226, 1).Value = 5
$ws.Cells.Item(226, 2).Value = 16
$ws.Cells.Item(226, 3).Value = 6
$ws.Cells.Item(226, 4).Value = 4
$ws.Cells.Item(226, 5).Formula = "=B226+D226"

$ws.Cells.Item(227, 1).Value = 2
$ws.Cells.Item(227, 2).Value = 13
$ws.Cells.Item(227, 3).Value = 3
$ws.Cells.Item(227, 4).Value = 7
$ws.Cells.Item(227, 5).Formula = "=B227+D227"

$ws.Cells.Item(228, 1).Value = 4
$ws.Cells.Item(228, 2).Value = 12
$ws.Cells.Item(228, 3).Value = 3
$ws.Cells.Item(228, 4).Value = 8
$ws.Cells.Item(228, 5).Formula = "=B228+D228"

$ws.Cells.Item(229, 1).Value = 8
$ws.Cells.Item(229, 2).Value = 14
$ws.Cells.Item(229, 3).Value = 6
$ws.Cells.Item(229, 4).Value = 6
$ws.Cells.Item(229, 5).Formula = "=B229+D229"

$ws.Cells.Item(230, 1).Value = 5
$ws.Cells.Item(230, 2).Value = 13
$ws.Cells.Item(230, 3).Value = 4
$ws.Cells.Item(230, 4).Value = 7
$ws.Cells.Item(230, 5).Formula = "=B230+D230"

$ws.Cells.Item(231, 1).Value = 5
$ws.Cells.Item(231, 2).Value = 15
$ws.Cells.Item(231, 3).Value = 4
$ws.Cells.Item(231, 4).Value = 5
$ws.Cells.Item(231, 5).Formula = "=B231+D231"

$ws.Cells.Item(232, 1).Value = 4
$ws.Cells.Item(232, 2).Value = 12
$ws.Cells.Item(232, 3).Value = 5
$ws.Cells.Item(232, 4).Value = 8
$ws.Cells.Item(232, 5).Formula = "=B232+D232"

$ws.Cells.Item(233, 1).Value = 3
$ws.Cells.Item(233, 2).Value = 15
$ws.Cells.Item(233, 3).Value = 5
$ws.Cells.Item(233, 4).Value = 5
$ws.Cells.Item(233, 5).Formula = "=B233+D233"

$ws.Cells.Item(234, 1).Value = 3
$ws.Cells.Item(234, 2).Value = 8
$ws.Cells.Item(234, 3).Value = 4
$ws.Cells.Item(234, 4).Value = 12
$ws.Cells.Item(234, 5).Formula = "=B234+D234"

$ws.Cells.Item(235, 1).Value = 3
$ws.Cells.Item(235, 2).Value = 16
$ws.Cells.Item(235, 3).Value = 4
$ws.Cells.Item(235, 4).Value = 4
$ws.Cells.Item(235, 5).Formula = "=B235+D235"

$ws.Cells.Item(236, 1).Value = 5
$ws.Cells.Item(236, 2).Value = 4
$ws.Cells.Item(236, 3).Value = 6
$ws.Cells.Item(236, 4).Value = 16
$ws.Cells.Item(236, 5).Formula = "=B236+D236"

$ws.Cells.Item(237, 1).Value = 7
$ws.Cells.Item(237, 2).Value = 14
$ws.Cells.Item(237, 3).Value = 4
$ws.Cells.Item(237, 4).Value = 6
$ws.Cells.Item(237, 5).Formula = "=B237+D237"

$ws.Cells.Item(238, 1).Value = 5
$ws.Cells.Item(238, 2).Value = 5
$ws.Cells.Item(238, 3).Value = 3
$ws.Cells.Item(238, 4).Value = 15
$ws.Cells.Item(238, 5).Formula = "=B238+D238"

$ws.Cells.Item(239, 1).Value = 3
$ws.Cells.Item(239, 2).Value = 8
$ws.Cells.Item(239, 3).Value = 4
$ws.Cells.Item(239, 4).Value = 12
$ws.Cells.Item(239, 5).Formula = "=B239+D239"

$ws.Cells.Item(240, 1).Value = 3
$ws.Cells.Item(240, 2).Value = 12
$ws.Cells.Item(240, 3).Value = 4
$ws.Cells.Item(240, 4).Value = 8
$ws.Cells.Item(240, 5).Formula = "=B240+D240"

$ws.Cells.Item(241, 1).Value = 5
$ws.Cells.Item(241, 2).Value = 13
$ws.Cells.Item(241, 3).Value = 4
$ws.Cells.Item(241, 4).Value = 7
$ws.Cells.Item(241, 5).Formula = "=B241+D241"

$ws.Cells.Item(242, 1).Value = 6
$ws.Cells.Item(242, 2).Value = 18
$ws.Cells.Item(242, 3).Value = 5
$ws.Cells.Item(242, 4).Value = 2
$ws.Cells.Item(242, 5).Formula = "=B242+D242"

$ws.Cells.Item(243, 1).Value = 3
$ws.Cells.Item(243, 2).Value = 7
$ws.Cells.Item(243, 3).Value = 4
$ws.Cells.Item(243, 4).Value = 13
$ws.Cells.Item(243, 5).Formula = "=B243+D243"

$ws.Cells.Item(244, 1).Value = 3
$ws.Cells.Item(244, 2).Value = 13
$ws.Cells.Item(244, 3).Value = 4
$ws.Cells.Item(244, 4).Value = 7
$ws.Cells.Item(244, 5).Formula = "=B244+D244"

$ws.Cells.Item(245, 1).Value = 2
$ws.Cells.Item(245, 2).Value = 8
$ws.Cells.Item(245, 3).Value = 3
$ws.Cells.Item(245, 4).Value = 12
$ws.Cells.Item(245, 5).Formula = "=B245+D245"

$ws.Cells.Item(246, 1).Value = 9
$ws.Cells.Item(246, 2).Value = 15
$ws.Cells.Item(246, 3).Value = 4
$ws.Cells.Item(246, 4).Value = 5
$ws.Cells.Item(246, 5).Formula = "=B246+D246"

$ws.Range("A247").Select()

